# Auto-generated edit script applying the cryptos.xlsx data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.207.09'
$ws.Range("E2").Value = '  +1.81%  '
$ws.Range("D3").Value = '3.386.08'
$ws.Range("E3").Value = '  +1.40%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '586.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '180.07'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.59%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.79%  '
$ws.Range("E9").Value = '  +6.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '48.62'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.62%  '
$ws.Range("E12").Value = '  +2.90%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '675.92'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.82%  '
$ws.Range("E14").Value = '  +1.96%  '
$ws.Range("D15").Value = '3.932.16'
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = '69.253.80'
$ws.Range("E16").Value = '  +1.92%  '
$ws.Range("D17").Value = '3.392.32'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("E18").Value = '  +1.78%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '17.66'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.36%  '
$ws.Range("E20").Value = '  +2.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.900'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.26%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.08'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.36%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '103.51'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.51%  '
$ws.Range("E25").Value = '  +0.59%  '
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("E27").Value = '  +0.56%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '34.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.54%  '
$ws.Range("E30").Value = '  -1.41%  '
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("B32").Value = 'dogwifhat'
$ws.Range("C32").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.68'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +11.62%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '553.61'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.23%  '
$ws.Range("E34").Value = '  +0.51%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '57.95'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.98%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("D37").Value = '3.689.46'
$ws.Range("E37").Value = '  -0.66%  '
$ws.Range("E38").Value = '  +6.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.96%  '
$ws.Range("E40").Value = '  +1.32%  '
$ws.Range("D41").Value = '0.0₃0703'
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.66'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.21%  '
$ws.Range("E43").Value = '  +0.43%  '
$ws.Range("E44").Value = '  +3.93%  '
$ws.Range("E45").Value = '  -2.23%  '
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.129'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.52%  '
$ws.Range("E48").Value = '  +5.43%  '
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '131.84'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.59%  '
$ws.Range("E51").Value = '  -1.16%  '
